# Updated cryptos list on Mon May 15 10:42:11 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" '27.639.66'
Set-TextCell $ws "E2" '  -0.64%  '
Set-TextCell $ws "D3" '1.846.14'
Set-TextCell $ws "E3" '  -1.10%  '
Set-TextCell $ws "D4" '1.010'
Set-TextCell $ws "E4" '  -2.89%  '
Set-TextCell $ws "D5" '318.68'
Set-TextCell $ws "E5" '  -1.96%  '
Set-TextCell $ws "D6" '1.010'
Set-TextCell $ws "E6" '  -2.60%  '
Set-TextCell $ws "D7" '0.4305'
Set-TextCell $ws "E7" '  -2.76%  '
Set-TextCell $ws "D8" '0.3748'
Set-TextCell $ws "E8" '  -1.46%  '
Set-TextCell $ws "D9" '0.07345'
Set-TextCell $ws "E9" '  -1.77%  '
Set-TextCell $ws "D10" '0.8796'
Set-TextCell $ws "E10" '  -0.77%  '
Set-TextCell $ws "D11" '21.57'
Set-TextCell $ws "E11" '  -1.00%  '
Set-TextCell $ws "D12" '1.856.19'
Set-TextCell $ws "E12" '  -0.92%  '
Set-TextCell $ws "D13" '6.726'
Set-TextCell $ws "E13" '  -0.65%  '
Set-TextCell $ws "D14" '5.449'
Set-TextCell $ws "E14" '  -2.07%  '
Set-TextCell $ws "D15" '0.07111'
Set-TextCell $ws "E15" '  -1.61%  '
Set-TextCell $ws "D16" '87.61'
Set-TextCell $ws "E16" '  +4.43%  '
Set-TextCell $ws "E17" '  -2.57%  '
Set-TextCell $ws "D18" '0.000008975'
Set-TextCell $ws "E18" '  -2.22%  '
Set-TextCell $ws "D19" '1.010'
Set-TextCell $ws "E19" '  -2.61%  '
Set-TextCell $ws "D20" '15.47'
Set-TextCell $ws "E20" '  -0.60%  '
Set-TextCell $ws "D21" '27.673.56'
Set-TextCell $ws "E21" '  -0.59%  '
Set-TextCell $ws "D22" '5.255'
Set-TextCell $ws "E22" '  -1.34%  '
Set-TextCell $ws "D23" '11.15'
Set-TextCell $ws "E23" '  -1.94%  '
Set-TextCell $ws "D24" '2.081.22'
Set-TextCell $ws "E24" '  -1.27%  '
Set-TextCell $ws "D25" '2.036'
Set-TextCell $ws "E25" '  +2.07%  '
Set-TextCell $ws "D26" '155.45'
Set-TextCell $ws "E26" '  -2.09%  '
Set-TextCell $ws "D27" '18.54'
Set-TextCell $ws "E27" '  -1.87%  '
Set-TextCell $ws "D28" '2.132'
Set-TextCell $ws "E28" '  +7.10%  '
Set-TextCell $ws "D29" '5.379'
Set-TextCell $ws "E29" '  +0.57%  '
Set-TextCell $ws "D30" '120.39'
Set-TextCell $ws "E30" '  +2.07%  '
Set-TextCell $ws "D31" '0.08919'
Set-TextCell $ws "E31" '  -1.90%  '
Set-TextCell $ws "D32" '1.226'
Set-TextCell $ws "E32" '  +0.72%  '
Set-TextCell $ws "D33" '0.7765'
Set-TextCell $ws "E33" '  -0.40%  '
Set-TextCell $ws "D34" '4.553'
Set-TextCell $ws "E34" '  -0.56%  '
Set-TextCell $ws "D35" '2.909'
Set-TextCell $ws "E35" '  -6.58%  '
Set-TextCell $ws "B36" 'TrustWalletToken'
Set-TextCell $ws "C36" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws "D36" '1.140'
Set-TextCell $ws "E36" '  -1.47%  '
Set-TextCell $ws "B37" 'Frax'
Set-TextCell $ws "C37" 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell $ws "D37" '1.010'
Set-TextCell $ws "E37" '  -2.68%  '
Set-TextCell $ws "D38" '0.05332'
Set-TextCell $ws "E38" '  -0.43%  '
Set-TextCell $ws "D39" '0.01970'
Set-TextCell $ws "E39" '  -1.49%  '
Set-TextCell $ws "D40" '7.213'
Set-TextCell $ws "E40" '  +4.07%  '
Set-TextCell $ws "D41" '2.864'
Set-TextCell $ws "E41" '  -0.34%  '
Set-TextCell $ws "D42" '0.5162'
Set-TextCell $ws "E42" '  -0.86%  '
Set-TextCell $ws "D43" '0.1676'
Set-TextCell $ws "E43" '  -1.32%  '
Set-TextCell $ws "D44" '8.918'
Set-TextCell $ws "E44" '  +2.62%  '
Set-TextCell $ws "D45" '110.46'
Set-TextCell $ws "E45" '  +0.51%  '
Set-TextCell $ws "D46" '10.70'
Set-TextCell $ws "E46" '  -0.12%  '
Set-TextCell $ws "D47" '0.4726'
Set-TextCell $ws "E47" '  +0.18%  '
Set-TextCell $ws "B48" 'Cronos'
Set-TextCell $ws "C48" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws "D48" '0.06495'
Set-TextCell $ws "E48" '  +0.38%  '
Set-TextCell $ws "B49" 'NEARProtocol'
Set-TextCell $ws "C49" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws "D49" '1.699'
Set-TextCell $ws "E49" '  -1.76%  '
Set-TextCell $ws "D50" '1.010'
Set-TextCell $ws "E50" '  -2.87%  '
Set-TextCell $ws "E51" '  -0.97%  '
